# Apply the "Solve Leetcode - 295. Find Median from Data Stream - 2 Heaps" edit
# to the Journal.xlsx / Neetcode 150 tracker workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neetcode 150")

# ---------------------------------------------------------------------------
# 1) Colour-code the DIFFICULTY column (B) using Excel's built-in cell
#    styles: Easy -> Good (green), Medium / Medium (!!!) -> Neutral (yellow),
#    Hard -> Bad (red). This matches the new fonts/fills/cellStyle "Bad"
#    entry that shows up in styles.xml.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 16; $r++) {
    $diff = $ws.Cells.Item($r, 2).Value2
    if ($diff -eq "Easy") {
        $ws.Cells.Item($r, 2).Style = "Good"
    } elseif ($diff -eq "Hard") {
        $ws.Cells.Item($r, 2).Style = "Bad"
    } else {
        # "Medium" and "Medium (!!!)"
        $ws.Cells.Item($r, 2).Style = "Neutral"
    }
}

# ---------------------------------------------------------------------------
# 2) Tweak the notes for "297. Serialize and Deserialize Binary Tree"
#    (row 10) - just a trailing period was added to the explanation.
# ---------------------------------------------------------------------------
$d10 = @"
First serialize into an array using preorder traversal and then return it as string with ",".join(serialized). Next split the serialized over "," and iterate over values (maybe store the vals into an iterator with vals = iter(data) and rebuild the same way you serialized using preorder.
"@
$ws.Range("D10").Value = $d10
$ws.Rows.Item(10).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 3) Add the new row (17) for "295. Find Median from Data Stream".
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "Heap/PQ"

$ws.Range("B17").Value = "Hard"
$ws.Range("B17").Style = "Bad"

$ws.Range("C17").Value = "295. Find Median from Data Stream"
$ws.Range("C17").Style = "Neutral"
$ws.Hyperlinks.Add($ws.Range("C17"), "https://leetcode.com/problems/find-median-from-data-stream/", "", "", "295. Find Median from Data Stream")

$d17 = @"
We are essentially gonna divide the list into two sorted halves in non decreasing order with two heaps - 
1) Left Heap will be a max heap because we need to check the rightmost element for mantaining the order or calculating median
2) Right Heap will be a min heap because we need to check its leftmost element for aforementioned reasons
With leftHeap[0] being <= rightHeap[0], and the heaps will be approximately the same size. And we'll balance the heaps whenever we add an element.
First check if the heaps mantain the leftHeap[-1] being <= rightHeap[0] order, if not the pop fromleft and push to right in a while loop.
Then, If the size of heaps differs by more than one element then move the top from min/max heap to the other.
Lastly, for calculating median, if lenghts of heaps is not equal then there are odd number of elements so return that from left or right heap, otherwise there are even elements so return (leftHeap[0] + rightHeap[0]) / 2
"@
$ws.Range("D17").Value = $d17
$ws.Rows.Item(17).RowHeight = 129.6

# ---------------------------------------------------------------------------
# 4) Update the sheet view / selection to match (scrolled down to show the
#    new row, with D17 selected) and fix up the used-range dimension.
# ---------------------------------------------------------------------------
$ws.Range("D17").Select()
$excel.ActiveWindow.ScrollRow = 14
